# Fix the numbered-list style used in the Id3 algorithm steps on the
# "Improved Id3 algorithm" slide: the auto-numbering scheme changes from
# "1." (arabicPeriod) to "1)" (arabicParenR) for every numbered step.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)

# "Content Placeholder 2" holds the numbered algorithm steps.
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

$count = $tr.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $bullet = $para.ParagraphFormat.Bullet
    if ($bullet.Type -eq 2) {
        # ppBulletNumbered — switch "arabicPeriod" (1.) to "arabicParenR" (1))
        $bullet.Style = 2
    }
}
